# Update cryptos list values to reflect refreshed market data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.440.41"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3
$ws.Range("D3").Value = "'3.290.23"
$ws.Range("E3").Value = "  +0.00%  "

# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.38%  "

# Row 5
$ws.Range("D5").Value = "'577.96"
$ws.Range("E5").Value = "  +4.24%  "

# Row 6
$ws.Range("D6").Value = "'182.78"
$ws.Range("E6").Value = "  -2.61%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").Value = "'3.286.31"
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("E9").Value = "  -2.38%  "

# Row 10
$ws.Range("D10").Value = "'0.175"
$ws.Range("E10").Value = "  -4.90%  "

# Row 11
$ws.Range("D11").Value = "'0.569"
$ws.Range("E11").Value = "  -2.77%  "

# Row 12
$ws.Range("D12").Value = "'46.40"
$ws.Range("E12").Value = "  -1.87%  "

# Row 13
$ws.Range("D13").Value = "'0.0000263"
$ws.Range("E13").Value = "  -2.40%  "

# Row 14
$ws.Range("D14").Value = "'632.58"
$ws.Range("E14").Value = "  +3.01%  "

# Row 15
$ws.Range("D15").Value = "'3.815.46"
$ws.Range("E15").Value = "  -0.31%  "

# Row 16
$ws.Range("E16").Value = "  -2.53%  "

# Row 17
$ws.Range("D17").Value = "'65.569.29"

# Row 18
$ws.Range("E18").Value = "  +0.14%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.65"
$ws.Range("E19").Value = "  -1.90%  "

# Row 20
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "'3.287.81"
$ws.Range("E20").Value = "  -0.43%  "

# Row 21
$ws.Range("D21").Value = "'10.92"
$ws.Range("E21").Value = "  -0.12%  "

# Row 22
$ws.Range("E22").Value = "  -2.21%  "

# Row 23
$ws.Range("D23").Value = "'17.94"
$ws.Range("E23").Value = "  -2.47%  "

# Row 24
$ws.Range("D24").Value = "'100.31"
$ws.Range("E24").Value = "  -1.58%  "

# Row 25
$ws.Range("D25").Value = "'4.96"
$ws.Range("E25").Value = "  +0.25%  "

# Row 26
$ws.Range("D26").Value = "'3.95"
$ws.Range("E26").Value = "  +0.60%  "

# Row 27
$ws.Range("D27").Value = "'2.73"
$ws.Range("E27").Value = "  +0.03%  "

# Row 28
$ws.Range("D28").Value = "'9.34"
$ws.Range("E28").Value = "  -2.34%  "

# Row 29
$ws.Range("D29").Value = "'30.63"
$ws.Range("E29").Value = "  +1.62%  "

# Row 30
$ws.Range("E30").Value = "  -3.38%  "

# Row 31
$ws.Range("D31").Value = "'6.47"
$ws.Range("E31").Value = "  -0.17%  "

# Row 32
$ws.Range("D32").Value = "'574.81"
$ws.Range("E32").Value = "  +3.49%  "

# Row 33
$ws.Range("D33").Value = "'3.67"
$ws.Range("E33").Value = "  -8.98%  "

# Row 34
$ws.Range("E34").Value = "  -1.79%  "

# Row 35
$ws.Range("D35").Value = "'3.845.77"
$ws.Range("E35").Value = "  +0.72%  "

# Row 36
$ws.Range("D36").Value = "'0.104"
$ws.Range("E36").Value = "  -1.16%  "

# Row 37
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("E38").Value = "  -3.09%  "

# Row 39
$ws.Range("E39").Value = "  -2.49%  "

# Row 40
$ws.Range("B40").Value = "ApeXProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  +4.97%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'32.34"
$ws.Range("E41").Value = "  -4.47%  "

# Row 42
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "'0.0₃0678"
$ws.Range("E42").Value = "  -5.92%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'3.11"
$ws.Range("E43").Value = "  -5.80%  "

# Row 44
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.58"
$ws.Range("E44").Value = "  -4.52%  "

# Row 45
$ws.Range("D45").Value = "'0.331"
$ws.Range("E45").Value = "  -1.39%  "

# Row 46
$ws.Range("E46").Value = "  -3.74%  "

# Row 47
$ws.Range("D47").Value = "'3.03"
$ws.Range("E47").Value = "  -4.78%  "

# Row 48
$ws.Range("E48").Value = "  +0.29%  "

# Row 49
$ws.Range("E49").Value = "  -1.73%  "

# Row 50
$ws.Range("D50").Value = "'2.51"
$ws.Range("E50").Value = "  -2.26%  "

# Row 51
$ws.Range("D51").Value = "'129.34"
$ws.Range("E51").Value = "  +5.54%  "
